# Daily attendance processing - reorder names in 'Recorded By' (column G) values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud'
$ws.Range("G3").Value = 'Dr. Gehan Adel, Dr. Manar Montaser, Administrator, Dr. Alshimaa Atef'
$ws.Range("G4").Value = 'Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Hanan Ragab, Dr. Majorelle Magdy, Dr. Asmaa Reda'
$ws.Range("G7").Value = 'Dr. Amal Awwad, Dr. Safa Hany'
$ws.Range("G9").Value = 'Dr. Yasmeena Fattoh, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna'
$ws.Range("G12").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G13").Value = 'Dr. Youstina Gamil, Dr. Sarah Mahdy'
$ws.Range("G17").Value = 'Dr. Marian Samir, Dr. Enas Omran, Dr. Walaa Ghanima'
$ws.Range("G18").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Eman Samir Gabry'
$ws.Range("G19").Value = 'Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Yasmin'
$ws.Range("G20").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G21").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud'
$ws.Range("G22").Value = 'Dr. Gehan Adel, Dr. Manar Montaser, Administrator, Dr. Alshimaa Atef'
$ws.Range("G23").Value = 'Dr. Nourhan Mahmoud, Dr. Shimaa Ahmad Mekki, Dr. Menna tuâ€™Allah Medhat, Dr. Heba Mahmoud Ali, Dr. Hanan Ragab, Dr. Majorelle Magdy, Dr. Asmaa Reda'
$ws.Range("G26").Value = 'Dr. Amal Awwad, Dr. Safa Hany'
$ws.Range("G28").Value = 'Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Marwa Mustafa, Dr. Sarah Abdelmohsen, Dr. Basma Hamed'
$ws.Range("G29").Value = 'Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G31").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G36").Value = 'Dr. Marian Samir, Dr. Enas Omran, Dr. Walaa Ghanima'
$ws.Range("G37").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Eman Samir Gabry'
$ws.Range("G38").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G39").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G40").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Hend Mahmoud'
$ws.Range("G41").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef, Dr. Hend Mahmoud'
$ws.Range("G42").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G43").Value = 'Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Kerelos Zareef'
$ws.Range("G45").Value = 'Dr. Amal Awwad, Dr. Safa Hany'
$ws.Range("G47").Value = 'Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Merna Said, Dr. Amira Ibrahim'
$ws.Range("G48").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Fatma Shoukry, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Sarah Abdelmohsen, Dr. Merna Said'
$ws.Range("G49").Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Range("G50").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G54").Value = 'Dr. Afaf Abdallah, Dr. Amr Saeed'
$ws.Range("G56").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Eman Samir Gabry'
$ws.Range("G57").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G58").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G59").Value = 'Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Nesma, Dr. Servinaz Sayed Mohammad'
$ws.Range("G60").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef, Dr. Hend Mahmoud'
$ws.Range("G61").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Majorelle Magdy, Dr. Asmaa Reda'
$ws.Range("G63").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G66").Value = 'Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G67").Value = 'Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G75").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Eman Samir Gabry'
$ws.Range("G76").Value = 'Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Yasmin'
$ws.Range("G77").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G78").Value = 'Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Nesma, Dr. Servinaz Sayed Mohammad'
$ws.Range("G79").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef, Dr. Hend Mahmoud'
$ws.Range("G80").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Nahla Nagiub, Dr. Majorelle Magdy, Dr. Asmaa Reda'
$ws.Range("G83").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G85").Value = 'Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Marina Youhanna, Dr. Amira Ibrahim'
$ws.Range("G86").Value = 'Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Fatma Shoukry, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Sarah Abdelmohsen, Dr. Merna Said'
$ws.Range("G88").Value = 'Dr. Mona Ibrahim Hussein, Dr. Dalia Tarek Elwan, Dr. Heba Al-Sayed Mohammad'
$ws.Range("G89").Value = 'Dr. Youstina Gamil, Dr. Sarah Mahdy'
$ws.Range("G94").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Eman Samir Gabry'
$ws.Range("G95").Value = 'Dr. Marina Sorial, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Yasmin'
$ws.Range("G96").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
$ws.Range("G97").Value = 'Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Mohammad El-Tanany, Dr. Asmaa Reda, Dr. Nesma, Dr. Servinaz Sayed Mohammad'
$ws.Range("G98").Value = 'Dr. Amira Sobhy, Dr. Shimaa Ahmad Mekki, Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Alshimaa Atef, Dr. Hend Mahmoud'
$ws.Range("G99").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range("G100").Value = 'Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Nada Mohammad, Dr. Menna tu''Alllah Mohammad, Dr. Abeer Ragab, Dr. Kerelos Zareef'
$ws.Range("G101").Value = 'Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany'
$ws.Range("G102").Value = 'Dr. Amal Awwad, Dr. Safa Hany'
$ws.Range("G104").Value = 'Dr. Arwa Al-Sayed, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Maryam Ahmad, Dr. Merna Said, Dr. Amira Ibrahim'
$ws.Range("G113").Value = 'Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Eman Samir Gabry'
$ws.Range("G115").Value = 'Dr. Yasmin, Dr. Marina Sorial, Dr. Neveen Nashaat, Dr. Nardine, Dr. Monica, Dr. Marina Atef, Dr. Remon'
